$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(3770, 4685, 4685, 4699, 4777, 4777, 4777, 4887, 4929, 4934, 4989, 5064, 5306, 5306)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 3).Value = $v
    $row++
}
